# Applies two content edits to the "Problem Solving" document:
#
#  1. Paragraph 13 ("The subgoal is to basically have it so ...") :
#     "subgoal" -> "sub-goal", written out as three separate runs
#     ("The sub" / "-" / "goal is to ...") with the spell-check
#     proofErr wrapper removed.
#
#  2. The lone bookmark paragraph that follows the socks "Problem:"
#     paragraph is expanded into three paragraphs:
#       - a blank paragraph
#       - "The constraints for this are that it is dark and you
#          can't see the color."
#       - "The sub-goal is to have one pair " <bookmark _GoBack> "of
#          all three colors."
#
# Note: this engine coalesces adjacent same-formatted runs that are
# inserted back-to-back, just like real Word normally would merge
# freshly typed text. Real authoring sessions end up with separate
# <w:r> elements in these spots anyway (e.g. because of an
# intervening spell-check pass or a bookmark), so we use a
# momentary bookmark as a run-splitting "wedge": insert the text, plant
# a bookmark at the seam, then delete the bookmark once the seam is
# baked in as a run boundary.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "subgoal" -> "sub-goal" (paragraph 13), split into 3 runs
# ---------------------------------------------------------------------

$p13 = $d.Paragraphs(13).Range

# Replace the paragraph's text (this also drops the proofErr spell-check
# tags, since they belonged to the run(s) being removed).
$textOnly13 = $d.Range($p13.Start, $p13.End - 1)
$textOnly13.Delete()

$p13 = $d.Paragraphs(13).Range
$ins13 = $d.Range($p13.Start, $p13.Start)
$ins13.InsertAfter("The subgoal is to basically have it so that the cat and the parrot are not left alone on the same side at the same time. Also the parrot and the seed are not left alone.")

# Split "The sub" | "-" | "goal is to ..." into separate runs using a
# pair of temporary bookmarks as seam markers, then remove them.
$p13 = $d.Paragraphs(13).Range
$seamPos = $p13.Start + 7   # right after "The sub"

$d.Bookmarks.Add("TmpSeamA", $d.Range($seamPos, $seamPos)) | Out-Null
$d.Range($seamPos, $seamPos).InsertAfter("-")
$d.Bookmarks.Add("TmpSeamB", $d.Range($seamPos, $seamPos)) | Out-Null

$d.Bookmarks("TmpSeamA").Delete()
$d.Bookmarks("TmpSeamB").Delete()

# ---------------------------------------------------------------------
# 2) Socks section: list the constraints + sub-goal
# ---------------------------------------------------------------------

# Paragraph 33 is the paragraph that currently holds only the _GoBack
# bookmark (right after the socks "Problem:" paragraph). Insert two new
# paragraphs ahead of it.
$p33 = $d.Paragraphs(33).Range
$p33.InsertParagraphBefore()
$p33.InsertParagraphBefore()

# The first of the two new paragraphs (now index 33) should stay blank.
# InsertParagraphBefore leaves behind an empty run; clean it out so the
# paragraph has no runs at all.
$p33 = $d.Paragraphs(33).Range
$insBlank = $d.Range($p33.Start, $p33.Start)
$insBlank.InsertAfter("X")
$p33 = $d.Paragraphs(33).Range
$d.Range($p33.Start, $p33.End - 1).Delete()

# The second new paragraph (now index 34) gets the constraints sentence.
$p34 = $d.Paragraphs(34).Range
$ins34 = $d.Range($p34.Start, $p34.Start)
$ins34.InsertAfter("The constraints for this are that it is dark and you can’t see the color.")

# The third paragraph (now index 35) is the original bookmark paragraph.
# Remove the bookmark, type the full sentence, then re-plant the
# bookmark right where it belongs - between "pair " and "of".
$d.Bookmarks("_GoBack").Delete()

$p35 = $d.Paragraphs(35).Range
$ins35 = $d.Range($p35.Start, $p35.Start)
$ins35.InsertAfter("The sub-goal is to have one pair of all three colors.")

$p35 = $d.Paragraphs(35).Range
$bmPos = $p35.Start + 33   # length of "The sub-goal is to have one pair "
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos)) | Out-Null
